$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 788, shifting rows 788:829 down to 789:830
# (matches <dimension ref="A1:D829"/> -> <dimension ref="A1:D830"/>).
$ws.Rows.Item(788).Insert()

# Column A holds date-like text (e.g. "2026/12/29") stored verbatim as a
# string, not a real date serial. Force the cell to Text format first so
# assigning "2026/02/06" doesn't get auto-converted into a date value, then
# reset the style back to Normal so no stray number-format/style is left
# attached to the cell (matching the plain, unstyled sibling data cells).
$ws.Cells.Item(788, 1).NumberFormat = "@"
$ws.Cells.Item(788, 1).Value = "2026/02/06"
$ws.Cells.Item(788, 1).Style = "Normal"

$ws.Cells.Item(788, 2).Value = "金"
$ws.Cells.Item(788, 3).Value = 3
$ws.Cells.Item(788, 4).Value = 201
